# Rename the two worksheets to prefix them with "V4.4" (tabs were renamed for
# the V4.4 Shenhe -> ShenheFrostFlower remap draft). Renaming the first sheet
# also auto-updates the _xlnm._FilterDatabase defined name, which refers to
# it by name.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$ws1.Name = "V4.4Shenhe to ShenheFrostFlower"
$ws2.Name = "V4.4SheheFrostFlower to Shenhe"

# Scroll sheet 1 down (the pane was frozen on row 1; the visible top row
# moves from row 17 to row 68) before switching focus to sheet 2, which
# becomes the active/selected tab.
$ws1.Activate()
$ws1.Range("A68").Select()

$ws2.Activate()
$ws2.Range("D21").Select()
